{"js": "// Apply the \"Iteration 2/3 sample sprints\" deadline-wording updates.\n//\n// Each change below replaces the text of one paragraph/run with the\n// updated wording from the commit. We use Range.search() scoped to the\n// document body (matchCase + wholeWords off, since we are matching\n// punctuation-bearing substrings) and then Range.insertText(\u2026, \"Replace\")\n// to rewrite just the matched span, leaving surrounding runs/paragraph\n// formatting untouched.\n\nasync function replaceOnce(body, searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// Week 9 bullets: \"N\" -> \"N/N*\"\nawait replaceOnce(\n  body,\n  \"Population of N relationship and with-clause information\",\n  \"Population of N/N* relationship and with-clause information\"\n);\n\nawait replaceOnce(\n  body,\n  \"Evaluation of queries with multiple clauses, focusing on N relationship, with-clauses\",\n  \"Evaluation of queries with multiple clauses, focusing on N/N* relationship, with-clauses\"\n);\n\n// Week 10 intro line: drop the leading \"N*/\" before \"A/A*\"\nawait replaceOnce(\n  body,\n  \"/QP] Help out with test cases, or N*/A/A* generation, or optimisation\",\n  \"/QP] Help out with test cases, or A/A* generation, or optimisation\"\n);\n\n// Week 10 bullets: \"N*/A\" -> \"A/A*\"\nawait replaceOnce(\n  body,\n  \"Population of N*/A relationship\",\n  \"Population of A/A* relationship\"\n);\n\nawait replaceOnce(\n  body,\n  \"Evaluation of queries with multiple clauses, focusing on N*/A relationship, tuple\",\n  \"Evaluation of queries with multiple clauses, focusing on A/A* relationship, tuple\"\n);\n\n// Week 11 bullets: \"A*\" -> \"NB/NB*\" / \"base\"\nawait replaceOnce(\n  body,\n  \"Population of A* relationship\",\n  \"Population of NB/NB* relationship\"\n);\n\nawait replaceOnce(\n  body,\n  \"Evaluation of queries with optimisation, focusing on A* correctness and finding out which optimisation strategy works best\",\n  \"Evaluation of queries with optimisation, focusing on base correctness and finding out which optimisation strategy works best\"\n);\n\n// Week 12 bullet: \"NB/NB*\" -> \"AB/AB*\"\nawait replaceOnce(\n  body,\n  \"Population of NB/NB* relationship extension\",\n  \"Population of AB/AB* relationship extension\"\n);\n\n// Week 13 no longer has a \"Features to be implemented\u2026\" feature list \u2014\n// remove the bold heading paragraph plus its two bullet paragraphs.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst targetAnchor = \"[G] Conduct presentation rehearsal after code + report submission\";\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetAnchor) {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find Week 13 anchor paragraph\");\n}\n\n// The next three paragraphs are the \"Features to be implemented\u2026\",\n// \"Population of AB/AB* relationship extension\" and \"Evaluation of\n// queries with optimisation, focusing on extension correctness\" \u2014\n// delete them (back-to-front so indices stay valid).\nfor (let i = anchorIndex + 3; i >= anchorIndex + 1; i--) {\n  items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Apply the \"Iteration 2/3 sample sprints\" deadline-wording updates\n# using the Word COM object model (Find/Replace + paragraph deletion).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($searchText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceOne = 1 (MatchCase=True, Forward=True, Wrap=wdFindContinue, Replace=wdReplaceOne)\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null\n}\n\n# Week 9 bullets: \"N\" -> \"N/N*\"\nReplace-Text \"Population of N relationship and with-clause information\" `\n    \"Population of N/N* relationship and with-clause information\"\n\nReplace-Text \"Evaluation of queries with multiple clauses, focusing on N relationship, with-clauses\" `\n    \"Evaluation of queries with multiple clauses, focusing on N/N* relationship, with-clauses\"\n\n# Week 10 intro line: drop the leading \"N*/\" before \"A/A*\"\nReplace-Text \"/QP] Help out with test cases, or N*/A/A* generation, or optimisation\" `\n    \"/QP] Help out with test cases, or A/A* generation, or optimisation\"\n\n# Week 10 bullets: \"N*/A\" -> \"A/A*\"\nReplace-Text \"Population of N*/A relationship\" `\n    \"Population of A/A* relationship\"\n\nReplace-Text \"Evaluation of queries with multiple clauses, focusing on N*/A relationship, tuple\" `\n    \"Evaluation of queries with multiple clauses, focusing on A/A* relationship, tuple\"\n\n# Week 11 bullets: \"A*\" -> \"NB/NB*\" / \"base\"\nReplace-Text \"Population of A* relationship\" `\n    \"Population of NB/NB* relationship\"\n\nReplace-Text \"Evaluation of queries with optimisation, focusing on A* correctness and finding out which optimisation strategy works best\" `\n    \"Evaluation of queries with optimisation, focusing on base correctness and finding out which optimisation strategy works best\"\n\n# Week 12 bullet: \"NB/NB*\" -> \"AB/AB*\"\nReplace-Text \"Population of NB/NB* relationship extension\" `\n    \"Population of AB/AB* relationship extension\"\n\n# Week 13 no longer has a \"Features to be implemented\u2026\" feature list \u2014\n# remove the bold heading paragraph plus its two bullet paragraphs that\n# followed \"[G] Conduct presentation rehearsal after code + report submission\".\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"[G] Conduct presentation rehearsal after code + report submission\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find Week 13 anchor paragraph\"\n}\n\nfor ($i = $anchorIndex + 3; $i -ge $anchorIndex + 1; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
